$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "About": rebuild the sources/notes block with the new BNEF
# record-low article, a real hyperlink, and a single recalibrated
# USD-normalisation factor (replacing the old two-row 1.29 / (1/0.951)
# pair with one literal value used as a multiplier).
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.UsedRange.Clear() | Out-Null

$about.Columns("A").ColumnWidth = 18.71

$about.Range("A1").Value = "BPP Battery Pack Price"
$about.Range("A1").Font.Bold = $true

$about.Range("A3").Value = "Sources:"
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = "BNEF"

$about.Range("B4").Value = "Lithium-ion Battery Pack Prices Hit Record Low of `$139/kWh"

$about.Range("B5").Value = 2023
$about.Range("B5").HorizontalAlignment = -4131

$about.Range("B6").Value = "https://about.bnef.com/blog/lithium-ion-battery-pack-prices-hit-record-low-of-139-kwh/#:~:text=Given%20this%2C%20BNEF%20expects%20average,and%20%2480%2FkWh%20in%202030."
$about.Hyperlinks.Add($about.Range("B6"), "https://about.bnef.com/blog/lithium-ion-battery-pack-prices-hit-record-low-of-139-kwh/", ":~:text=Given%20this%2C%20BNEF%20expects%20average,and%20%2480%2FkWh%20in%202030.") | Out-Null

$about.Range("A8").Value = "Notes:"
$about.Range("A8").Font.Bold = $true

$about.Range("A9").Value = "The EPS applies endogenous learning for battery pack prices in years where the battery pack price is listed as 0."

$about.Range("A11").Value = "2023 to 2012"
$about.Range("B11").Value = 0.75350342301658668
$about.Range("B11").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# Sheet "BPP": 2021/2022/2023 battery pack prices recalibrated against
# the new BNEF figures ($150, $161, $139) using the new About!$B$11
# multiplier instead of the old division by About!$A$16 / $A$17.
# ---------------------------------------------------------------------
$bpp = $wb.Worksheets.Item("BPP")

$bpp.Range("B2").Formula = '=150*About!$B$11'
$bpp.Range("C2").Formula = '=161*About!$B$11'
$bpp.Range("D2").Formula = '=139*About!$B$11'

$bpp.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "SYBPP": 2020 battery pack price recalibrated the same way.
# ---------------------------------------------------------------------
$sybpp = $wb.Worksheets.Item("SYBPP")

$sybpp.Range("B2").Formula = '=160*About!B11'

$sybpp.Range("D13").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore "About" as the active/selected sheet + cell last, so its
# tabSelected flag and final selection match the saved workbook state.
# ---------------------------------------------------------------------
$about.Activate() | Out-Null
$about.Range("B21").Select() | Out-Null
